$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-31 with re-sorted data
$ws.Range("A2").Value = 45
$ws.Range("B2").Value = "T931383615020"
$ws.Range("A3").Value = 106
$ws.Range("B3").Value = "V931101109012"
$ws.Range("A4").Value = 23
$ws.Range("B4").Value = "N931252508050"
$ws.Range("A5").Value = 30
$ws.Range("B5").Value = "V931414517045"
$ws.Range("A6").Value = 0
$ws.Range("B6").Value = "L394201008038"
$ws.Range("A7").Value = 64
$ws.Range("B7").Value = "K931252910051"
$ws.Range("A8").Value = 127
$ws.Range("B8").Value = "V931101109041"
$ws.Range("A9").Value = 74
$ws.Range("B9").Value = "V333218013124"
$ws.Range("A10").Value = 77
$ws.Range("B10").Value = "L931383612038"
$ws.Range("A11").Value = 82
$ws.Range("B11").Value = "R931101109037"
$ws.Range("A12").Value = 90
$ws.Range("B12").Value = "D931252710019"
$ws.Range("A13").Value = 104
$ws.Range("B13").Value = "V931100509030"
$ws.Range("A14").Value = 107
$ws.Range("B14").Value = "U931253114004"
$ws.Range("A15").Value = 114
$ws.Range("B15").Value = "F931100609012"
$ws.Range("A16").Value = 123
$ws.Range("B16").Value = "U931325208066"
$ws.Range("A17").Value = 63
$ws.Range("B17").Value = "T887690719015"
$ws.Range("A18").Value = 60
$ws.Range("B18").Value = "P931101109055"
$ws.Range("A19").Value = 65
$ws.Range("B19").Value = "J931101109071"
$ws.Range("A20").Value = 50
$ws.Range("B20").Value = "V931321008075"
$ws.Range("A21").Value = 2
$ws.Range("B21").Value = "Y931252508056"
$ws.Range("A22").Value = 36
$ws.Range("B22").Value = "A380404117005"
$ws.Range("A23").Value = 33
$ws.Range("B23").Value = "K931100609063"
$ws.Range("A24").Value = 24
$ws.Range("B24").Value = "F931252108061"
$ws.Range("A25").Value = 6
$ws.Range("B25").Value = "M931325212046"
$ws.Range("A26").Value = 21
$ws.Range("B26").Value = "U931258914007"
$ws.Range("A27").Value = 7
$ws.Range("B27").Value = "X931235209022"
$ws.Range("A28").Value = 18
$ws.Range("B28").Value = "M931412019009"
$ws.Range("A29").Value = 15
$ws.Range("B29").Value = "C931101008023"
$ws.Range("A30").Value = 58
$ws.Range("B30").Value = "Q931253109015"
$ws.Range("A31").Value = 122
$ws.Range("B31").Value = "F931252509025"

# Add new rows 32-46, copying formatting from row 31 (A column) for consistent style
$ws.Range("A31").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A32").Value = 120
$ws.Range("B32").Value = "U931252110024"
$ws.Range("A31").Copy()
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A33").Value = 118
$ws.Range("B33").Value = "B931235209044"
$ws.Range("A31").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("A34").Value = 113
$ws.Range("B34").Value = "Q931252108040"
$ws.Range("A31").Copy()
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("A35").Value = 111
$ws.Range("B35").Value = "J931384210007"
$ws.Range("A31").Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A36").Value = 101
$ws.Range("B36").Value = "A931252108046"
$ws.Range("A31").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A37").Value = 53
$ws.Range("B37").Value = "Z931325309007"
$ws.Range("A31").Copy()
$ws.Range("A38").PasteSpecial(-4122)
$ws.Range("A38").Value = 95
$ws.Range("B38").Value = "K931325309035"
$ws.Range("A31").Copy()
$ws.Range("A39").PasteSpecial(-4122)
$ws.Range("A39").Value = 94
$ws.Range("B39").Value = "Y931325210014"
$ws.Range("A31").Copy()
$ws.Range("A40").PasteSpecial(-4122)
$ws.Range("A40").Value = 84
$ws.Range("B40").Value = "H931321309010"
$ws.Range("A31").Copy()
$ws.Range("A41").PasteSpecial(-4122)
$ws.Range("A41").Value = 34
$ws.Range("B41").Value = "M931259308029"
$ws.Range("A31").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("A42").Value = 73
$ws.Range("B42").Value = "R931253116053"
$ws.Range("A31").Copy()
$ws.Range("A43").PasteSpecial(-4122)
$ws.Range("A43").Value = 39
$ws.Range("B43").Value = "N931253409013"
$ws.Range("A31").Copy()
$ws.Range("A44").PasteSpecial(-4122)
$ws.Range("A44").Value = 47
$ws.Range("B44").Value = "K931252509020"
$ws.Range("A31").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$ws.Range("A45").Value = 105
$ws.Range("B45").Value = "K931383410019"
$ws.Range("A31").Copy()
$ws.Range("A46").PasteSpecial(-4122)
$ws.Range("A46").Value = 128
$ws.Range("B46").Value = "B931400418001"
